# The document has two logos (the Pearson Edexcel logo in the footers and
# the BTec logo in the headers), each appearing once in the "first page"
# header/footer and once in the "default" (primary) header/footer. This
# edit simply swaps the `Name` each inline picture reports/stores in the
# drawing's docPr (the BTec logo was "image1.jpg" and becomes "image2.jpg";
# the Pearson logo was "image2.png" and becomes "image1.png") - the
# pictures themselves (and their relationship targets) are unchanged.

$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2

foreach ($sec in $d.Sections) {

    # --- Headers: BTec logo, image1.jpg -> image2.jpg -------------------
    foreach ($hfIndex in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage)) {
        $hdr = $sec.Headers.Item($hfIndex)
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    # --- Footers: Pearson logo, image2.png -> image1.png ----------------
    foreach ($hfIndex in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage)) {
        $ftr = $sec.Footers.Item($hfIndex)
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}

Write-Output "done"
